$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Refresh the panel_query_time timestamps on the data sheet (rows 2-64, column F)
$queryTimes = @(
    "2021-10-05 14:34:22.456100",
    "2021-10-05 14:34:22.456108",
    "2021-10-05 14:34:22.456111",
    "2021-10-05 14:34:22.456114",
    "2021-10-05 14:34:22.456117",
    "2021-10-05 14:34:22.456120",
    "2021-10-05 14:34:22.456122",
    "2021-10-05 14:34:22.456125",
    "2021-10-05 14:34:22.456128",
    "2021-10-05 14:34:22.456130",
    "2021-10-05 14:34:22.456133",
    "2021-10-05 14:34:22.456135",
    "2021-10-05 14:34:22.456138",
    "2021-10-05 14:34:22.456140",
    "2021-10-05 14:34:22.456143",
    "2021-10-05 14:34:22.456145",
    "2021-10-05 14:34:22.456148",
    "2021-10-05 14:34:22.456151",
    "2021-10-05 14:34:22.456153",
    "2021-10-05 14:34:22.456156",
    "2021-10-05 14:34:22.456159",
    "2021-10-05 14:34:22.456162",
    "2021-10-05 14:34:22.456164",
    "2021-10-05 14:34:22.456167",
    "2021-10-05 14:34:22.456170",
    "2021-10-05 14:34:22.456172",
    "2021-10-05 14:34:22.456175",
    "2021-10-05 14:34:22.456177",
    "2021-10-05 14:34:22.456180",
    "2021-10-05 14:34:22.456182",
    "2021-10-05 14:34:22.456185",
    "2021-10-05 14:34:22.456187",
    "2021-10-05 14:34:22.456190",
    "2021-10-05 14:34:22.456193",
    "2021-10-05 14:34:22.456196",
    "2021-10-05 14:34:22.456198",
    "2021-10-05 14:34:22.456201",
    "2021-10-05 14:34:22.456203",
    "2021-10-05 14:34:22.456206",
    "2021-10-05 14:34:22.456208",
    "2021-10-05 14:34:22.456211",
    "2021-10-05 14:34:22.456214",
    "2021-10-05 14:34:22.456216",
    "2021-10-05 14:34:22.456219",
    "2021-10-05 14:34:22.456221",
    "2021-10-05 14:34:22.456224",
    "2021-10-05 14:34:22.456227",
    "2021-10-05 14:34:22.456229",
    "2021-10-05 14:34:22.456232",
    "2021-10-05 14:34:22.456234",
    "2021-10-05 14:34:22.456237",
    "2021-10-05 14:34:22.456239",
    "2021-10-05 14:34:22.456242",
    "2021-10-05 14:34:22.456245",
    "2021-10-05 14:34:22.456248",
    "2021-10-05 14:34:22.456250",
    "2021-10-05 14:34:22.456253",
    "2021-10-05 14:34:22.456255",
    "2021-10-05 14:34:22.456258",
    "2021-10-05 14:34:22.456260",
    "2021-10-05 14:34:22.456263",
    "2021-10-05 14:34:22.456266",
    "2021-10-05 14:34:22.456268"
)
for ($i = 0; $i -lt $queryTimes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $queryTimes[$i]
}


# Add the new "metadata" sheet (placed after "data", matching the panelapp export format)
$meta = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Limb Girdle Muscular Dystrophy"
$meta.Range("C2").Value = 3071
$meta.Range("D2").Value = "0.59"
$meta.Range("E2").Value = "2021-06-19T02:11:21.021874Z"
$meta.Range("F2").Value = "2021-10-05 14:34:22.452505"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3071/?format=json"

# Mirror the bordered/bold/centered header style already used on the "data" sheet
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore "data" as the active sheet/selection
$ws.Activate()

Write-Output "done"
